# Add season-record columns (Wins, Losses, Ties) to the stats sheet.
#
# The sheet currently spans A1:AC45. We extend it to A1:AF45 by adding
# three new columns: AD = Wins, AE = Losses, AF = Ties.
# Every data row (2-45) gets the same season record: 91 wins, 71 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting used by the other header cells (bold, centered,
# bordered) by copying the style from an existing header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows ----------------------------------------------------------
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 71
    $ws.Cells.Item($r, 32).Value = 0
}

$ws.Range("A1").Select()
